# Commit: "mise en commentaires des parametres de recherche qu'on n'utilise
# plus" (update-mapping-pn13).
#
# On the "Metadata" sheet:
#   - refresh the "Date" property's value
#   - insert a new "Jurisdiction" property (with an empty value) right
#     after "Contact" and before "Description", pushing every row below
#     it down by one (Description/Purpose/Copyright/Immutable).
#
# The "Include from Medications" sheet is untouched; it merely reflects
# the shared-string renumbering caused by the new rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" property value (row 8, column B) -------------------
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# --- Make room for the new "Jurisdiction" row at row 11 --------------------
# Shift rows 11..14 (Description, Purpose, Copyright, Immutable) down to
# rows 12..15, working bottom-up so we never clobber a row before reading
# it. Formatting is carried along explicitly (copy + paste-format) so the
# shifted rows keep the sheet's existing cell style instead of picking up
# a brand new one.
for ($r = 14; $r -ge 11; $r--) {
    $destRow = $r + 1
    $aVal = $ws.Cells.Item($r, 1).Value()
    $bVal = $ws.Cells.Item($r, 2).Value()

    $ws.Cells.Item($destRow, 1).Value = $aVal
    $ws.Cells.Item($destRow, 2).Value = $bVal

    $ws.Range("A" + $r + ":B" + $r).Copy()
    $ws.Range("A" + $destRow + ":B" + $destRow).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# --- Populate the freed-up row 11 with the new "Jurisdiction" property -----
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
